$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily record was inserted for 2026/02/21 (Saturday) between the
# existing row 853 (2026/02/21, 16:00, ...) and the old row 854
# (2026/12/29). Inserting a whole row shifts every row below it down by
# one, which is exactly what the diff shows (old row 854 data now lives
# in row 855, old 855 -> 856, ... old 895 -> 896), and bumps the used
# range from A1:D895 to A1:D896.
$ws.Rows.Item(854).Insert()

$ws.Range("A854").Value = "'2026/02/21"
$ws.Range("B854").Value = "土"
$ws.Range("C854").Value = 20
$ws.Range("D854").Value = 201
